$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column C width (matches xl/worksheets/sheet1.xml <cols> entry, ~23.66 chars)
$ws.Columns.Item(3).ColumnWidth = 22.8

# Row 115: section header
$ws.Cells.Item(115, 1).Value = "Table 4 -- Recombinants"

# Row 116
$ws.Cells.Item(116, 1).Value = "RF2k/1b"
$ws.Cells.Item(116, 2).Value = 3186
$ws.Cells.Item(116, 3).Value = "AY587845"

# Row 117
$ws.Cells.Item(117, 1).Value = "RF2i/6p"
$ws.Cells.Item(117, 2).Value = "3405-3464"
$ws.Cells.Item(117, 3).Value = "DQ155560"

# Row 118
$ws.Cells.Item(118, 1).Value = "RF2b/1b_1"
$ws.Cells.Item(118, 2).Value = 3456
$ws.Cells.Item(118, 3).Value = "DQ364460"

# Row 119
$ws.Cells.Item(119, 1).Value = "RF2/5"
$ws.Cells.Item(119, 2).Value = "3366-3389"
$ws.Cells.Item(119, 3).Value = "AM408911"

# Row 120
$ws.Cells.Item(120, 1).Value = "RF2b/6w"
$ws.Cells.Item(120, 2).Value = 3429
$ws.Cells.Item(120, 3).Value = "EU643835"

# Row 121
$ws.Cells.Item(121, 1).Value = "RF2b/1b_2"
$ws.Cells.Item(121, 2).Value = 3432
$ws.Cells.Item(121, 3).Value = "AB622121"

# Row 122
$ws.Cells.Item(122, 1).Value = "RF2b/1a"
$ws.Cells.Item(122, 2).Value = "3429-3440"
$ws.Cells.Item(122, 3).Value = "JF779679"

# Row 123
$ws.Cells.Item(123, 1).Value = "RF2b/1b_3"
$ws.Cells.Item(123, 2).Value = "3286-3293"
$ws.Cells.Item(123, 3).Value = "AB677530"

# Row 124
$ws.Cells.Item(124, 1).Value = "RF2b/1b_4"
$ws.Cells.Item(124, 2).Value = "3286-3293"
$ws.Cells.Item(124, 3).Value = "AB677527"

# Apply the same cell style (s="2", Helvetica font) used by the rest of the table,
# touching only the cells that actually hold values (matches the diff exactly).
$srcStyle = $ws.Cells.Item(113, 1)
$targets = @(
    @(115,1),
    @(116,1), @(116,2), @(116,3),
    @(117,1), @(117,2), @(117,3),
    @(118,1), @(118,2), @(118,3),
    @(119,1), @(119,2), @(119,3),
    @(120,1), @(120,2), @(120,3),
    @(121,1), @(121,2), @(121,3),
    @(122,1), @(122,2), @(122,3),
    @(123,1), @(123,2), @(123,3),
    @(124,1), @(124,2), @(124,3)
)
foreach ($t in $targets) {
    $cell = $ws.Cells.Item($t[0], $t[1])
    $cell.Font.Name = $srcStyle.Font.Name
    $cell.Font.Size = $srcStyle.Font.Size
}

# Update the view: scroll + selection to match
$ws.Application.ActiveWindow.ScrollRow = 97
$ws.Range("A116:C124").Select()
